# Remove the "Mfg Part #" column (E) and the "Type" column (originally H,
# now G after the first deletion) from the BOM table, then remove the now
# -empty trailing rows (4-8) that were left over as filler/formatting rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column E ("Mfg Part #") - shifts F..I left by one.
$ws.Range("E1").EntireColumn.Delete()

# Delete entire column G ("Type", originally H before the shift above)
$ws.Range("G1").EntireColumn.Delete()

# Remove the now-unused filler rows 4-8 (rows that had no real content).
$ws.Range("A4:H8").EntireRow.Delete()

# Update the view: scroll to column C and select F10 (per saved view state).
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("F10").Select()
